$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Formatting first (xlPasteFormats = -4122): stamp each new/changed
# row with the same cell styles used by the existing "group start"
# (style 4/5, border-less) and "group end" (style 7/8, bottom border)
# rows elsewhere in the table.
# ------------------------------------------------------------------

# Row 38 becomes a "group end" row (border) once column A gets a value.
$ws.Range("A34:E34").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial(-4122) | Out-Null

# Row 39: new "group start" row.
$ws.Range("A35:E35").Copy() | Out-Null
$ws.Range("A39:E39").PasteSpecial(-4122) | Out-Null

# Row 40: new "group end" row, only column A populated.
$ws.Range("A34:E34").Copy() | Out-Null
$ws.Range("A40:E40").PasteSpecial(-4122) | Out-Null

# Row 41: new "group start" row.
$ws.Range("A35:E35").Copy() | Out-Null
$ws.Range("A41:E41").PasteSpecial(-4122) | Out-Null

# Row 42: new row (same style family as 41, column A populated).
$ws.Range("A35:E35").Copy() | Out-Null
$ws.Range("A42:E42").PasteSpecial(-4122) | Out-Null

# Rows 43-44: continuation rows with no value in column A at all
# (so only paste formats into B:E, leaving A43/A44 nonexistent).
$ws.Range("B35:E35").Copy() | Out-Null
$ws.Range("B43:E43").PasteSpecial(-4122) | Out-Null
$ws.Range("B35:E35").Copy() | Out-Null
$ws.Range("B44:E44").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# Cell values, entered in the same order the source workbook's
# shared-string table shows them (keeps the sheet content identical
# cell-for-cell to the target edit).
# ------------------------------------------------------------------
$ws.Range("A38").Value = "SCRIPT/T01P02A/um2506.ssb"

$ws.Range("C39").Value = ' Good luck on the graduation\nexam. Hey, hey, hey!'
$ws.Range("A39").Value = "SCRIPT/P01P01A/us0102.ssb"
$ws.Range("D39").Value = ' Удачи вам на выпускном экзамене.\nЭй, эй, эй!'
$ws.Range("E39").Value = ' Ôäàœé âàí îà âúðôòëîïí üëèàíåîå.\nÜê, üê, üê!'

$ws.Range("A40").Value = "SCRIPT/P01P01A/us0104.ssb"

$ws.Range("C41").Value = ' Visiting the guild, hey, hey?!'
$ws.Range("C42").Value = ' Makes me glad to see you!\nHey, hey!'
$ws.Range("C43").Value = ' Even after you graduate, you\''re\nfree to look up jobs on the boards and go after\noutlaws.'
$ws.Range("C44").Value = ' So keep taking on those jobs,\nhey, hey.'

$ws.Range("A41").Value = "SCRIPT/G01P03A/us0111.ssb"

$ws.Range("D41").Value = ' Посещаете гильдию, эй, эй?!'
$ws.Range("D42").Value = ' Я рад вас здесь видеть! Эй, эй!'
$ws.Range("D43").Value = ' Даже после выпуска вы по\nпрежнему можете смотреть задания на доске\nобъявлений и ловить негодяев.'
$ws.Range("D44").Value = ' Поэтому, беритесь за всё, что\nвам угодно, эй, эй.'

$ws.Range("E41").Value = ' Ðïòåþàåóå ãéìûäéý, üê, üê?!'
$ws.Range("E42").Value = ' Ÿ ñàä âàò èäåòû âéäåóû! Üê, üê!'
$ws.Range("E43").Value = ' Äàçå ðïòìå âúðôòëà âú ðï\nðñåçîåíô íïçåóå òíïóñåóû èàäàîéÿ îà äïòëå\nïáùÿâìåîéê é ìïâéóû îåãïäÿåâ.'
$ws.Range("E44").Value = ' Ðïüóïíô, áåñéóåòû èà âòæ, œóï\nâàí ôãïäîï, üê, üê.'

$ws.Range("A42").Value = "SCRIPT/G01P03A/us3103.ssb"

# ------------------------------------------------------------------
# Numeric line-number column (not shared strings, order is cosmetic).
# ------------------------------------------------------------------
$ws.Range("B38").Value = 397
$ws.Range("B39").Value = 375
$ws.Range("B41").Value = 328
$ws.Range("B42").Value = 331
$ws.Range("B43").Value = 334
$ws.Range("B44").Value = 337

# ------------------------------------------------------------------
# Row heights to match the wrapped-text line counts of the new rows.
# ------------------------------------------------------------------
$ws.Rows.Item(38).RowHeight = 43.2
$ws.Rows.Item(39).RowHeight = 43.2
$ws.Rows.Item(40).RowHeight = 43.2
$ws.Rows.Item(41).RowHeight = 43.2
$ws.Rows.Item(42).RowHeight = 28.2
$ws.Rows.Item(43).RowHeight = 31.8
$ws.Rows.Item(44).RowHeight = 21.6

# ------------------------------------------------------------------
# Scroll / selection state, matching where the editor ended up.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("D43").Select() | Out-Null
